# Add "Bottom View of Binary Tree" row (row 12) to Sheet1, mirroring the
# pattern used by the existing "Top View of Binary Tree" row (row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B12: GFG/LC source -> "GFG"
$ws.Range("B12").Value = "GFG"

# C12: Question text -> new shared string
$ws.Range("C12").Value = "Bottom View of Binary Tree"

# D12: Language -> "Java/Python"
$ws.Range("D12").Value = "Java/Python"

# E12: Difficulty -> "Medium"
$ws.Range("E12").Value = "Medium"

# Copy formatting from row 11 cells so the new row matches the sheet's
# existing style (B: s=1, D: s=8, E: s=11, C: default style).
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("D11").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null

$ws.Range("E11").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update the active selection to C12, matching the diff.
$ws.Range("C12").Select() | Out-Null
